{"js": "// Replace the four body paragraphs of the letter and split the signature\n// line \"[Doctor Name]\" into \"[Name]\" / \"[Doctor]\" on two lines, per the\n// commit's updated letter copy.\n\nconst replacements = [\n  {\n    find: \"I am writing to you regarding [Patient Name], who is currently under my care. I have been monitoring [Patient Name]'s condition and have noticed a few concerning changes in their health.\",\n    replace: \"I am writing to you regarding [Patient Name], who is currently under my care. I wanted to provide you with an update on their condition and the treatment plan we have put in place.\"\n  },\n  {\n    find: \"[Patient Name] has been experiencing [list symptoms], which I believe are indicative of a more serious underlying condition. I am recommending that [Patient Name] undergo further testing to determine the cause of these symptoms.\",\n    replace: \"[Patient Name] has been diagnosed with [condition], and we are currently working to manage the symptoms and improve their overall health. We have prescribed [medication] and have been monitoring their progress closely.\"\n  },\n  {\n    find: \"I understand that this may be a difficult time for [Patient Name], and I want to assure you that I am here to provide support and guidance throughout the process. I am confident that with the right care and treatment, [Patient Name] will be able to make a full recovery.\",\n    replace: \"At this time, [Patient Name] is responding well to the treatment and is showing signs of improvement. We are continuing to monitor their progress and adjust the treatment plan as needed.\"\n  },\n  {\n    find: \"If you have any questions or concerns, please do not hesitate to contact me. I am available to discuss [Patient Name]'s condition and treatment options at any time.\",\n    replace: \"I understand that this can be a difficult time for [Patient Name] and their family. I want to assure you that we are doing everything we can to ensure their health and wellbeing. If you have any questions or concerns, please do not hesitate to contact me.\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + find);\n  }\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// \"[Doctor Name]\" -> \"[Name]\" + line break + \"[Doctor]\"\nconst sigResults = context.document.body.search(\"[Doctor Name]\", { matchCase: true });\nsigResults.load(\"items\");\nawait context.sync();\n\nif (sigResults.items.length === 0) {\n  throw new Error(\"Could not find text: [Doctor Name]\");\n}\nsigResults.items[0].insertText(\"[Name]\\u000b[Doctor]\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the four body paragraphs of the letter and split the signature\n# line \"[Doctor Name]\" into \"[Name]\" / \"[Doctor]\" on two lines, per the\n# commit's updated letter copy.\n\n$d = $word.ActiveDocument\n\nfunction Replace-LetterText($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-LetterText \"I am writing to you regarding [Patient Name], who is currently under my care. I have been monitoring [Patient Name]'s condition and have noticed a few concerning changes in their health.\" \"I am writing to you regarding [Patient Name], who is currently under my care. I wanted to provide you with an update on their condition and the treatment plan we have put in place.\"\n\nReplace-LetterText \"[Patient Name] has been experiencing [list symptoms], which I believe are indicative of a more serious underlying condition. I am recommending that [Patient Name] undergo further testing to determine the cause of these symptoms.\" \"[Patient Name] has been diagnosed with [condition], and we are currently working to manage the symptoms and improve their overall health. We have prescribed [medication] and have been monitoring their progress closely.\"\n\nReplace-LetterText \"I understand that this may be a difficult time for [Patient Name], and I want to assure you that I am here to provide support and guidance throughout the process. I am confident that with the right care and treatment, [Patient Name] will be able to make a full recovery.\" \"At this time, [Patient Name] is responding well to the treatment and is showing signs of improvement. We are continuing to monitor their progress and adjust the treatment plan as needed.\"\n\nReplace-LetterText \"If you have any questions or concerns, please do not hesitate to contact me. I am available to discuss [Patient Name]'s condition and treatment options at any time.\" \"I understand that this can be a difficult time for [Patient Name] and their family. I want to assure you that we are doing everything we can to ensure their health and wellbeing. If you have any questions or concerns, please do not hesitate to contact me.\"\n\nReplace-LetterText \"[Doctor Name]\" \"[Name]^l[Doctor]\"\n"}
